# Auto-generated edit script: updates Betfair back/lay odds cells
# for 2025-12-26 matches (rows 2-25) to match the new source snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2.04
$ws.Range("H2").Value = 3.55
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 4.3
$ws.Range("N2").Value = 6.4
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 2.8
$ws.Range("Q2").Value = 1.52
$ws.Range("R2").Value = 1.74
$ws.Range("S2").Value = 2.28
$ws.Range("T2").Value = 1.52
$ws.Range("U2").Value = 2.82
$ws.Range("V2").Value = 1.38
$ws.Range("W2").Value = 1.96
$ws.Range("X2").Value = 28
$ws.Range("Z2").Value = 32
$ws.Range("AA2").Value = 70
$ws.Range("AB2").Value = 16
$ws.Range("AD2").Value = 15.5
$ws.Range("AE2").Value = 34
$ws.Range("AI2").Value = 34
$ws.Range("AL2").Value = 24
$ws.Range("AN2").Value = 8.6
$ws.Range("AO2").Value = 21

# Row 3
$ws.Range("F3").Value = 1.01
$ws.Range("J3").Value = 1.2
$ws.Range("N3").Value = 1.1
$ws.Range("Q3").Value = 1.08
$ws.Range("V3").Value = 1.22

# Row 4
$ws.Range("F4").Value = 5.1
$ws.Range("G4").Value = 6.4
$ws.Range("H4").Value = 1.61
$ws.Range("I4").Value = 1.69
$ws.Range("N4").Value = 5.8

# Row 5
$ws.Range("I5").Value = 4.7
$ws.Range("L5").Value = 1.54
$ws.Range("N5").Value = 2.3
$ws.Range("Q5").Value = 2.88

# Row 6
$ws.Range("F6").Value = 1.43
$ws.Range("G6").Value = 1.5
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 4.7
$ws.Range("P6").Value = 1.6
$ws.Range("R6").Value = 1.24
$ws.Range("S6").Value = 4.1
$ws.Range("T6").Value = 2.36
$ws.Range("U6").Value = 1.57
$ws.Range("W6").Value = 3
$ws.Range("X6").Value = 1000
$ws.Range("AB6").Value = 990

# Row 7
$ws.Range("G7").Value = 1.16
$ws.Range("H7").Value = 27
$ws.Range("J7").Value = 1.2
$ws.Range("U7").Value = 1.45
$ws.Range("W7").Value = 3

# Row 8
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 8.4
$ws.Range("I8").Value = 12.5
$ws.Range("K8").Value = 6.6
$ws.Range("N8").Value = 6
$ws.Range("P8").Value = 2.68
$ws.Range("R8").Value = 1.68
$ws.Range("V8").Value = 1.09
$ws.Range("W8").Value = 3.6

# Row 9
$ws.Range("F9").Value = 3.25
$ws.Range("G9").Value = 3.65
$ws.Range("H9").Value = 2.16
$ws.Range("I9").Value = 2.38
$ws.Range("J9").Value = 3.4
$ws.Range("V9").Value = 1.72

# Row 10
$ws.Range("F10").Value = 2.34
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2.96
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 3.6
$ws.Range("K10").Value = 3.95
$ws.Range("U10").Value = 2.02
$ws.Range("W10").Value = 1.67

# Row 11
$ws.Range("F11").Value = 2.06
$ws.Range("G11").Value = 2.22
$ws.Range("K11").Value = 4.4
$ws.Range("L11").Value = 1.2
$ws.Range("R11").Value = 1.7
$ws.Range("S11").Value = 2.18
$ws.Range("T11").Value = 1.49
$ws.Range("V11").Value = 1.38
$ws.Range("W11").Value = 1.81
$ws.Range("AO11").Value = 21

# Row 12
$ws.Range("F12").Value = 2.36
$ws.Range("G12").Value = 2.58
$ws.Range("J12").Value = 3.15
$ws.Range("R12").Value = 1.28
$ws.Range("W12").Value = 1.64

# Row 13
$ws.Range("F13").Value = 6.8
$ws.Range("I13").Value = 1.51
$ws.Range("L13").Value = 1.22
$ws.Range("S13").Value = 2.3
$ws.Range("V13").Value = 2.96

# Row 14
$ws.Range("K14").Value = 4.1

# Row 15
$ws.Range("F15").Value = 1.2
$ws.Range("G15").Value = 1.25
$ws.Range("H15").Value = 11
$ws.Range("I15").Value = 15.5
$ws.Range("J15").Value = 7.6
$ws.Range("K15").Value = 10.5
$ws.Range("P15").Value = 3.3
$ws.Range("Q15").Value = 1.31
$ws.Range("R15").Value = 1.9
$ws.Range("T15").Value = 1.89
$ws.Range("U15").Value = 1.9
$ws.Range("V15").Value = 1.07
$ws.Range("AC15").Value = 23
$ws.Range("AF15").Value = 10.5
$ws.Range("AJ15").Value = 11
$ws.Range("AM15").Value = 190
$ws.Range("AN15").Value = 3.35

# Row 16
$ws.Range("W16").Value = 1.79
$ws.Range("Z16").Value = 30
$ws.Range("AG16").Value = 11.5

# Row 17
$ws.Range("F17").Value = 5.6
$ws.Range("H17").Value = 1.59
$ws.Range("I17").Value = 1.69
$ws.Range("J17").Value = 4.1
$ws.Range("K17").Value = 4.8

# Row 18
$ws.Range("F18").Value = 2.56
$ws.Range("G18").Value = 2.84
$ws.Range("I18").Value = 3.2
$ws.Range("K18").Value = 3.65
$ws.Range("L18").Value = 1.43
$ws.Range("R18").Value = 1.3
$ws.Range("S18").Value = 3.7
$ws.Range("V18").Value = 1.46
$ws.Range("W18").Value = 1.54

# Row 20
$ws.Range("I20").Value = 1.58
$ws.Range("M20").Value = 1.01
$ws.Range("R20").Value = 1.89
$ws.Range("T20").Value = 1.63

# Row 21
$ws.Range("H21").Value = 1.7
$ws.Range("K21").Value = 4.3

# Row 22
$ws.Range("G22").Value = 2.58
$ws.Range("I22").Value = 3.75
$ws.Range("J22").Value = 3.05
$ws.Range("K22").Value = 3.1
$ws.Range("Q22").Value = 2.5

# Row 23
$ws.Range("J23").Value = 3.8
$ws.Range("K23").Value = 4.3
$ws.Range("AB23").Value = 17
$ws.Range("AC23").Value = 11.5
$ws.Range("AK23").Value = 26
$ws.Range("AO23").Value = 25

# Row 24
$ws.Range("F24").Value = 1.19
$ws.Range("G24").Value = 1.23
$ws.Range("H24").Value = 14
$ws.Range("I24").Value = 18.5
$ws.Range("J24").Value = 8.4
$ws.Range("K24").Value = 12.5
$ws.Range("N24").Value = 10.5
$ws.Range("Q24").Value = 1.25
$ws.Range("S24").Value = 1.66
$ws.Range("T24").Value = 1.68
$ws.Range("U24").Value = 2.2
$ws.Range("W24").Value = 5.3
$ws.Range("Y24").Value = 110
$ws.Range("Z24").Value = 220
$ws.Range("AA24").Value = 1000
$ws.Range("AB24").Value = 19.5
$ws.Range("AC24").Value = 24
$ws.Range("AD24").Value = 1000
$ws.Range("AE24").Value = 210
$ws.Range("AH24").Value = 34
$ws.Range("AI24").Value = 150
$ws.Range("AJ24").Value = 12.5
$ws.Range("AL24").Value = 29
$ws.Range("AM24").Value = 140
$ws.Range("AN24").Value = 2.78

# Row 25
$ws.Range("J25").Value = 3.75
$ws.Range("K25").Value = 3.8
$ws.Range("P25").Value = 2.38
$ws.Range("R25").Value = 1.57
$ws.Range("T25").Value = 1.57
$ws.Range("U25").Value = 2.66
$ws.Range("AF25").Value = 20
$ws.Range("AI25").Value = 29
